# Fix error in reading GAMS parameters
#
# Two lines had errors when trying to read in GAMS parameters. The `damadj`
# term (not used in GAMS) is hard coded to 1 (already present in this sheet),
# and the `eqmat` term is added to the parameter table.
#
# The new "eqmat" parameter goes in its own row right after the "fco22x"
# row (row 30) and before the blank row that precedes the "CLIMATE DYNAMICS"
# section header -- i.e. it becomes the new row 31, pushing every
# subsequent row down by one (Excel automatically keeps all formula
# references, like the t2xco2-based B-column formula, pointing at the
# correct shifted cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 31, shifting "CLIMATE DYNAMICS" (and
# everything below it) down by one row.
$ws.Rows(31).Insert() | Out-Null

# Populate the new row with the "eqmat" parameter and its value.
$ws.Range("A31").Value = "eqmat"
$ws.Range("B31").Value = 588

# Match the author's final on-screen selection/viewport.
$ws.Range("A32").Select() | Out-Null
